## Fruta / hortaliza, semanal
## Inserts two new weekly price rows (874-875) for
## "Feria Lagunitas de Puerto Montt - Limón", pushing the existing
## historical rows (previously 874-963) down to 876-965.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 874, shifting everything
# below (874-963) down to (876-965).
$ws.Range("A874:T875").EntireRow.Insert()

# --- New row 874: 1a plateado -----------------------------------------
$ws.Range("A874").Value = 4
$ws.Range("B874").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C874").Value = "Los Lagos"
$ws.Range("D874").Value = 45212
$ws.Range("E874").Value = 10
$ws.Range("F874").Value = "Fruta"
$ws.Range("G874").Value = 100102
$ws.Range("H874").Value = "Cítricos"
$ws.Range("I874").Value = 100102003
$ws.Range("J874").Value = "Limón"
$ws.Range("K874").Value = "Sin especificar"
$ws.Range("L874").Value = "1a plateado"
$ws.Range("M874").Value = 1200
$ws.Range("N874").Value = 14000
$ws.Range("O874").Value = 15000
$ws.Range("P874").Value = 14500
$ws.Range("Q874").Value = "$/malla 18 kilos"
$ws.Range("R874").Value = "Provincia de Melipilla"
$ws.Range("S874").Value = 806
$ws.Range("T874").Value = 18

# --- New row 875: 2a plateado -----------------------------------------
$ws.Range("A875").Value = 4
$ws.Range("B875").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C875").Value = "Los Lagos"
$ws.Range("D875").Value = 45212
$ws.Range("E875").Value = 10
$ws.Range("F875").Value = "Fruta"
$ws.Range("G875").Value = 100102
$ws.Range("H875").Value = "Cítricos"
$ws.Range("I875").Value = 100102003
$ws.Range("J875").Value = "Limón"
$ws.Range("K875").Value = "Sin especificar"
$ws.Range("L875").Value = "2a plateado"
$ws.Range("M875").Value = 600
$ws.Range("N875").Value = 13000
$ws.Range("O875").Value = 13000
$ws.Range("P875").Value = 13000
$ws.Range("Q875").Value = "$/malla 18 kilos"
$ws.Range("R875").Value = "Provincia de Melipilla"
$ws.Range("S875").Value = 722
$ws.Range("T875").Value = 18
